# Edit script implementing:
#  1) Add a new bullet "Create CarSerializer class in serializers.py"
#     right before the "Makemigrations" bullet that follows "Create Car model".
#  2) Move <w:lastRenderedPageBreak/> from the "Test" run to the
#     "Add to path to Car's urls" run that precedes it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: insert "Create CarSerializer class in serializers.py" bullet
# ---------------------------------------------------------------------
# Find the "Makemigrations" paragraph that is immediately preceded by the
# "Create Car model" paragraph (there are two "Makemigrations" bullets in
# the document; this locates the second one). Track its 1-based index
# within $d.Paragraphs as we walk the collection.
$targetPos = -1
$pos = 0
$prevText = ""
foreach ($p in $d.Paragraphs) {
    $pos++
    $t = $p.Range.Text
    if (($t -match "Makemigrations") -and ($prevText -match "Create Car model")) {
        $targetPos = $pos
    }
    $prevText = $t
}

$target = $d.Paragraphs.Item($targetPos)

# Insert a new paragraph before it; InsertParagraphBefore() clones the
# paragraph formatting (style + numbering) of $target onto the new
# paragraph, matching the surrounding list items.
$target.Range.InsertParagraphBefore()

# After the insertion the brand-new (still empty) paragraph occupies the
# same document-wide index $targetPos used to hold $target; $target
# itself (and everything after it) shifted one slot later.
$newPara = $d.Paragraphs.Item($targetPos)
$newPara.Range.Text = "Create CarSerializer class in serializers.py"

# ---------------------------------------------------------------------
# Change 2: relocate <w:lastRenderedPageBreak/>
# ---------------------------------------------------------------------
# Locate the "Add to path to Car's urls" paragraph, which is immediately
# followed by a "Test" paragraph (there are several "Test" bullets in the
# doc; this pins down the right pair).
$addPara = $null
$testPara = $null
$prev = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($prev -ne $null) -and ($t -match "^Test\r?$") -and ($prev.Range.Text -match "Add to path to Car")) {
        $addPara = $prev
        $testPara = $p
    }
    $prev = $p
}

$fullRange = $d.Range($addPara.Range.Start, $testPara.Range.End)

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"
$apostrophe = [char]0x2019

$xml = "<w:p $ns w14:paraId=`"26F3E055`" w14:textId=`"47AB52AB`" w:rsidR=`"000110AD`" w:rsidRDefault=`"000110AD`" w:rsidP=`"000A23AD`">" +
       "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
       "<w:r><w:lastRenderedPageBreak/><w:t>Add to path to Car${apostrophe}s urls</w:t></w:r></w:p>" +
       "<w:p $ns w14:paraId=`"128F3654`" w14:textId=`"7066F6C6`" w:rsidR=`"000110AD`" w:rsidRDefault=`"000110AD`" w:rsidP=`"000A23AD`">" +
       "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
       "<w:r><w:t>Test</w:t></w:r></w:p>"

$fullRange.InsertXML($xml)
